# Rename the inline Pearson/BTEC logo pictures in every header/footer.
#
#   *.jpg  (BTec_Logo-Orange, in the headers)   image1.jpg -> image2.jpg
#   *.png  (PearsonLogo.png,  in the footers)   image2.png -> image1.png
#
# The pictures are identified by their (stable) alt-text / description
# rather than by header/footer index, so the script is robust no matter
# how wdHeaderFooterIndex happens to map onto the underlying parts.

$d = $word.ActiveDocument

function Rename-LogoInlineShapes($range) {
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $descr = $shape.AlternativeText

        if ($descr -eq "BTec_Logo-Orange") {
            $shape.Name = "image2.jpg"
        }
        elseif ($descr -like "*PearsonLogo.png") {
            $shape.Name = "image1.png"
        }
    }
}

foreach ($sec in $d.Sections) {
    for ($hfIndex = 1; $hfIndex -le 3; $hfIndex++) {
        $hdr = $sec.Headers.Item($hfIndex)
        if ($hdr.Exists) {
            Rename-LogoInlineShapes $hdr.Range
        }

        $ftr = $sec.Footers.Item($hfIndex)
        if ($ftr.Exists) {
            Rename-LogoInlineShapes $ftr.Range
        }
    }
}
